# Update VaR figures across the three report sheets.
# Commit message: "Removal of many unused functions" - the underlying
# calculation functions were pruned/changed, which altered the resulting
# VaR values that were pasted as static numbers into the report sheets.

$wb = $excel.ActiveWorkbook

# --- Total_VaR sheet ---
$wsTotal = $wb.Worksheets.Item("Total_VaR")
$wsTotal.Range("A2").Value = 619970.5599999964

# --- VaR_by_BUSINESS_LINE sheet ---
$wsBusiness = $wb.Worksheets.Item("VaR_by_BUSINESS_LINE")
$wsBusiness.Range("B2").Value = 374335.9999999985
$wsBusiness.Range("B3").Value = 701141.5200000077
$wsBusiness.Range("B4").Value = 25262

# --- VaR_by_METAL sheet ---
$wsMetal = $wb.Worksheets.Item("VaR_by_METAL")
$wsMetal.Range("B2").Value = 616426.0000000001
$wsMetal.Range("B3").Value = 25262
